# Apply the commit: "revision, added pyrolysis and additional figures"
#
# The "main" sheet gains a new boolean parameter row, "chemical_recycling_pyrolysis",
# inserted immediately after "chemical_recycling_gasification" (the existing row 9).
# Inserting the row pushes every row below it down by one (old row 10 "fossil_routes"
# becomes row 11, etc.), and the sheet's active selection ends up on B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Insert a new row at position 10, shifting existing rows 10-24 down to 11-25.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row with the new parameter.
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true

# Match the resulting selection/active cell shown in the saved workbook.
$ws.Activate() | Out-Null
$ws.Range("B9").Select() | Out-Null
